# Auto-generated Excel COM-interop script to apply cryptos list update
# (commit: "Updated cryptos list on Mon Jan 29 18:08:11 UTC 2024 with GitHub Actions")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.001.07"
$ws.Range("E2").Value = "  +1.99%  "
$ws.Range("D3").Value = "2.302.81"
$ws.Range("E3").Value = "  +1.48%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "309.34"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +1.16%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "100.34"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +4.46%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.534"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +0.81%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.508"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +3.59%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.94"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +1.92%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0817"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +2.88%  "
$ws.Range("E12").Value = "  +0.51%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.98"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +5.11%  "
$ws.Range("D14").Value = "2.658.92"
$ws.Range("E14").Value = "  +1.43%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.84"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +2.91%  "
$ws.Range("D16").Value = "2.302.28"
$ws.Range("E16").Value = "  +1.15%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.802"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +0.90%  "
$ws.Range("D18").Value = "42.963.46"
$ws.Range("E18").Value = "  +2.12%  "
$ws.Range("E19").Value = "  +0.39%  "
$ws.Range("D20").Value = "0.0₃0922"
$ws.Range("E20").Value = "  +1.62%  "
$ws.Range("E21").Value = "  +1.37%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "68.14"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +0.14%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "239.64"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +0.64%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.02"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +4.18%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.61"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +1.02%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.997"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.32%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.21"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +2.04%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "39.03"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +5.84%  "
$ws.Range("E29").Value = "  +1.11%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.12"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +0.56%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "169.56"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +5.94%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.33"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +1.44%  "
$ws.Range("E33").Value = "  +0.04%  "
$ws.Range("E34").Value = "  -1.93%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "17.71"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +3.69%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0739"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -0.30%  "
$ws.Range("E37").Value = "  +0.57%  "
$ws.Range("E38").Value = "  -0.05%  "
$ws.Range("E39").Value = "  +0.55%  "
$ws.Range("E40").Value = "  +1.06%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.20"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +4.52%  "
$ws.Range("B42").Value = "EnergySwap"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "19.62"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +3.16%  "
$ws.Range("B43").Value = "ApeXProtocol"
$ws.Range("C43").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.30"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -5.98%  "
$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0290"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +2.04%  "
$ws.Range("B45").Value = "Maker"
$ws.Range("C45").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D45").Value = "1.969.21"
$ws.Range("E45").Value = "  -0.95%  "
$ws.Range("E46").Value = "  +2.76%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.77"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -1.55%  "
$ws.Range("E48").Value = "  +19.15%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "55.06"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +3.61%  "
$ws.Range("B50").Value = "Stacks"
$ws.Range("C50").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.55"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +2.57%  "
$ws.Range("B51").Value = "RocketPoolETH"
$ws.Range("C51").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D51").Value = "2.525.65"
$ws.Range("E51").Value = "  +1.39%  "
